$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8733660130718954
$ws.Range("C2").Value = 0.9255411255411256
$ws.Range("D2").Value = 0.8986969314838169
$ws.Range("E2").Value = 1155

$ws.Range("B3").Value = 0.9556259904912837
$ws.Range("C3").Value = 0.9319938176197836
$ws.Range("D3").Value = 0.9436619718309859

$ws.Range("B4").Value = 0.8474148802017655
$ws.Range("C4").Value = 0.8538754764930114
$ws.Range("D4").Value = 0.850632911392405
$ws.Range("E4").Value = 787

$ws.Range("B5").Value = 0.8493150684931506
$ws.Range("C5").Value = 0.7065527065527065
$ws.Range("D5").Value = 0.7713841368584758
$ws.Range("E5").Value = 351

$ws.Range("B6").Value = 0.8816326530612245
$ws.Range("C6").Value = 0.8816326530612245
$ws.Range("D6").Value = 0.8816326530612245
$ws.Range("E6").Value = 0.8816326530612245

$ws.Range("B7").Value = 0.8814304880645237
$ws.Range("C7").Value = 0.8544907815516567
$ws.Range("D7").Value = 0.866093987891421

$ws.Range("B8").Value = 0.8816506328931242
$ws.Range("C8").Value = 0.8816326530612245
$ws.Range("D8").Value = 0.8805265935175524
